$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.803.95"
$ws.Range("E2").Value = "  +2.60%  "
$ws.Range("D3").Value = "2.529.87"
$ws.Range("E3").Value = "  -0.41%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "592.81"
$ws.Range("E5").Value = "  +2.29%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "176.76"
$ws.Range("E6").Value = "  +5.87%  "
$ws.Range("E7").Value = "  -0.09%  "
$ws.Range("E8").Value = "  +1.11%  "
$ws.Range("D9").Value = "2.528.31"
$ws.Range("E9").Value = "  -0.43%  "
$ws.Range("E10").Value = "  +1.64%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.164"
$ws.Range("E11").Value = "  +2.27%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.16"
$ws.Range("E12").Value = "  +0.72%  "
$ws.Range("E13").Value = "  -1.81%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "26.86"
$ws.Range("E14").Value = "  +1.59%  "
$ws.Range("D15").Value = "2.987.03"
$ws.Range("E15").Value = "  -0.63%  "
$ws.Range("E16").Value = "  +1.24%  "
$ws.Range("D17").Value = "67.636.36"
$ws.Range("E17").Value = "  +2.39%  "
$ws.Range("D18").Value = "2.513.33"
$ws.Range("E18").Value = "  -0.72%  "
$ws.Range("E19").Value = "  +5.04%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.46"
$ws.Range("E20").Value = "  +1.49%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "359.34"
$ws.Range("E21").Value = "  +3.84%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.20"
$ws.Range("E22").Value = "  +0.73%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.65"
$ws.Range("E23").Value = "  +2.46%  "
$ws.Range("E24").Value = "  +3.95%  "
$ws.Range("E26").Value = "  +4.29%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "70.87"
$ws.Range("E27").Value = "  +3.05%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.997"
$ws.Range("E28").Value = "  -0.18%  "
$ws.Range("D30").Value = "0.0₃0990"
$ws.Range("E30").Value = "  +1.58%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "555.49"
$ws.Range("E31").Value = "  +6.09%  "
$ws.Range("E32").Value = "  +1.07%  "
$ws.Range("E33").Value = "  +3.30%  "
$ws.Range("E34").Value = "  +2.90%  "
$ws.Range("E35").Value = "  +0.61%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.999"
$ws.Range("E36").Value = "  -0.07%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.48"
$ws.Range("E37").Value = "  +2.55%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "155.67"
$ws.Range("E38").Value = "  -0.65%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "18.76"
$ws.Range("E39").Value = "  +0.76%  "
$ws.Range("E40").Value = "  +1.78%  "
$ws.Range("B41").Value = "PolygonEcosystemToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.356"
$ws.Range("E41").Value = "  +0.75%  "
$ws.Range("B42").Value = "Stacks"
$ws.Range("C42").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.81"
$ws.Range("E42").Value = "  +3.36%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.18"
$ws.Range("E43").Value = "  +2.65%  "
$ws.Range("E44").Value = "  +5.81%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "147.52"
$ws.Range("E46").Value = "  +0.07%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.562"
$ws.Range("E47").Value = "  +1.11%  "
$ws.Range("D48").Value = "0.0₆0279"
$ws.Range("E48").Value = "  -0.84%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.72"
$ws.Range("E49").Value = "  +1.22%  "
$ws.Range("E50").Value = "  -0.16%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0757"
$ws.Range("E51").Value = "  +0.14%  "
